# Update the "Generate Report for Handback" timestamps.
$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date (row 2)
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 01:03:10"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 01:03:01"
$wsZhCn.Range("K2").Value = "2016-08-24 01:03:29"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 01:03:10"
$wsDeDe.Range("K2").Value = "2016-08-24 01:03:38"
